$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The stats for row 2 and row 3 (runs/balls/fours) were swapped:
#   Row 2: runs 22 -> 20, balls 13 -> 14, fours 2 -> 1
#   Row 3: runs 20 -> 22, balls 14 -> 13, fours 1 -> 2
# These columns are stored as text (numbers-as-text) in the sheet, so we
# write the new value through a TEXT() formula and then paste-special just
# the value back on top of itself. That keeps the cell a plain text value
# (matching the existing "t=str" cells) without leaving the cell in a
# formula state and without forcing a quote-prefixed / re-styled text cell.

function Set-TextNumber($cellAddr, $newValue) {
    $rng = $ws.Range($cellAddr)
    $rng.Formula = '=TEXT(' + $newValue + ',"0")'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

Set-TextNumber "C2" 20
Set-TextNumber "D2" 14
Set-TextNumber "E2" 1

Set-TextNumber "C3" 22
Set-TextNumber "D3" 13
Set-TextNumber "E3" 2
